# Generate Report for Handback
# - Flips the handback status text (cascades to Overview + both language sheets).
# - Fills in the "Latest Target File" / "Latest Handback File" columns (F/G) for
#   both language sheets, with hyperlinks matching the existing Source/Handoff
#   File link style.
# - Stamps the "Latest Handback DateTime" column (H) with real timestamps.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# zh-cn: Latest Handback DateTime
$zhcn.Range("H2").Value = "2016-03-24 08:49:45"
$zhcn.Range("H3").Value = "2016-03-24 08:49:45"

# de-de: Latest Handback DateTime
$dede.Range("H2").Value = "2016-03-24 08:49:52"
$dede.Range("H3").Value = "2016-03-24 08:49:52"

# zh-cn: Latest Target File (F) / Latest Handback File (G), rows 2 and 3
$zhcnMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/0a3a996dff101c17d296a87b4347f266e84e4dbf/e2e/34ded686-4006-40a0-a24a-57ef94237596.md"
$zhcnMd3Url = "https://github.com/OpenLocalizationTest/oltest/blob/0a3a996dff101c17d296a87b4347f266e84e4dbf/e2e/ffff11cb6873-c225-4731-8417-48de30a3441c.md"
$zhcnXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/72963e108308c1c39870dac2d05435e3b1434867/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/34ded686-4006-40a0-a24a-57ef94237596.5c6b38f4de7feb44bf2d40b2852ecd473571ec83.zh-cn.xlf"

$zhcn.Hyperlinks.Add($zhcn.Range("F2"), $zhcnMdUrl, "", "", "34ded686-4006-40a0-a24a-57ef94237596.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), $zhcnXlfUrl, "", "", "34ded686-4006-40a0-a24a-57ef94237596.5c6b38f4de7feb44bf2d40b2852ecd473571ec83.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), $zhcnMdUrl, "", "", "34ded686-4006-40a0-a24a-57ef94237596.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), $zhcnXlfUrl, "", "", "34ded686-4006-40a0-a24a-57ef94237596.5c6b38f4de7feb44bf2d40b2852ecd473571ec83.zh-cn.xlf")

# de-de: Latest Target File (F) / Latest Handback File (G), rows 2 and 3
$dedeMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/0a3a996dff101c17d296a87b4347f266e84e4dbf/e2e/34ded686-4006-40a0-a24a-57ef94237596.md"
$dedeXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c32d239dcec73e8605290148acaa033669ff389f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/34ded686-4006-40a0-a24a-57ef94237596.5c6b38f4de7feb44bf2d40b2852ecd473571ec83.de-de.xlf"

$dede.Hyperlinks.Add($dede.Range("F2"), $dedeMdUrl, "", "", "34ded686-4006-40a0-a24a-57ef94237596.md")
$dede.Hyperlinks.Add($dede.Range("G2"), $dedeXlfUrl, "", "", "34ded686-4006-40a0-a24a-57ef94237596.5c6b38f4de7feb44bf2d40b2852ecd473571ec83.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("F3"), $dedeMdUrl, "", "", "34ded686-4006-40a0-a24a-57ef94237596.md")
$dede.Hyperlinks.Add($dede.Range("G3"), $dedeXlfUrl, "", "", "34ded686-4006-40a0-a24a-57ef94237596.5c6b38f4de7feb44bf2d40b2852ecd473571ec83.de-de.xlf")
